$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "Provincia de Limari" (row 3: 44320-series, row 4: 44719-series)
# had their date / price values swapped between each other.
# Row 3 should now hold the values that were previously in row 4, and vice versa.

# Row 3 (becomes the 44719 / newer record)
$ws.Range("D3").Value = 44719
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20400
$ws.Range("S3").Value = 1133

# Row 4 (becomes the 44320 / older record)
$ws.Range("D4").Value = 44320
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 18800
$ws.Range("S4").Value = 1044
